$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("RunManager")

# --- Add two new test-case rows to the RunManager sheet ---
$ws2.Range("A5").Value = "validateLoginPageTitle"
$ws2.Range("B5").Value = "Yes"
$ws2.Range("A6").Value = "validateLoginWithIncorrectCredentials"
$ws2.Range("B6").Value = "No"

# New row 5 gets the same look as row 1 (wrapped, default font) -- copy formats only
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$ws2.Range("A1").Copy() | Out-Null
$ws2.Range("A5").PasteSpecial($xlPasteFormats) | Out-Null

# New row 6 gets the same look as the other test-case rows (A2:A4) -- copy formats only
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("A6").PasteSpecial($xlPasteFormats) | Out-Null

# Size the new rows to fit their (now wrapped) text
$ws2.Rows.Item(5).RowHeight = 17
$ws2.Rows.Item(6).RowHeight = 29

# RunManager is now the active/selected sheet and tab, with a new active selection
$ws2.Activate() | Out-Null
$ws2.Range("C10").Select() | Out-Null
